$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A4").Value = "referral_id"
$ws.Range("A17").Value = "genome_build"
